# Apply cryptos price/volume updates (Thu May 30 20:11:05 UTC 2024 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.300.16"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").Value = "3.733.71"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").Value = "'593.87"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").Value = "'166.98"
$ws.Range("E6").Value = "  -1.23%  "

$ws.Range("D7").Value = "3.733.82"
$ws.Range("E7").Value = "  -0.40%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").Value = "'0.519"
$ws.Range("E9").Value = "  -0.82%  "

$ws.Range("D10").Value = "'0.159"
$ws.Range("E10").Value = "  -3.44%  "

$ws.Range("D11").Value = "'6.45"
$ws.Range("E11").Value = "  -0.10%  "

$ws.Range("D12").Value = "'0.448"
$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("D13").Value = "'0.0000259"
$ws.Range("E13").Value = "  -5.78%  "

$ws.Range("D14").Value = "'36.26"
$ws.Range("E14").Value = "  -0.16%  "

$ws.Range("D15").Value = "4.367.10"
$ws.Range("E15").Value = "  -0.34%  "

$ws.Range("D16").Value = "3.747.71"
$ws.Range("E16").Value = "  -0.39%  "

$ws.Range("D17").Value = "68.375.19"
$ws.Range("E17").Value = "  +1.67%  "

$ws.Range("D18").Value = "'17.92"
$ws.Range("E18").Value = "  -3.14%  "

$ws.Range("D19").Value = "'7.00"
$ws.Range("E19").Value = "  -2.53%  "

$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("D21").Value = "'10.71"
$ws.Range("E21").Value = "  +2.22%  "

$ws.Range("D22").Value = "'466.58"
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").Value = "'0.700"
$ws.Range("E23").Value = "  -2.14%  "

$ws.Range("D24").Value = "'84.21"
$ws.Range("E24").Value = "  +0.60%  "

$ws.Range("D25").Value = "'0.0000144"
$ws.Range("E25").Value = "  -1.95%  "

$ws.Range("D26").Value = "'2.19"
$ws.Range("E26").Value = "  -0.31%  "

$ws.Range("D27").Value = "'12.06"
$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("D28").Value = "'10.10"
$ws.Range("E28").Value = "  -1.70%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").Value = "3.888.50"
$ws.Range("E30").Value = "  -0.41%  "

$ws.Range("D31").Value = "'2.79"
$ws.Range("E31").Value = "  -3.56%  "

$ws.Range("D32").Value = "'7.32"
$ws.Range("E32").Value = "  -4.02%  "

$ws.Range("D33").Value = "'29.90"
$ws.Range("E33").Value = "  -1.63%  "

$ws.Range("D34").Value = "'2.19"
$ws.Range("E34").Value = "  -1.47%  "

$ws.Range("D35").Value = "'9.27"
$ws.Range("E35").Value = "  +2.22%  "

$ws.Range("D37").Value = "3.695.09"
$ws.Range("E37").Value = "  -0.56%  "

$ws.Range("D38").Value = "'0.101"
$ws.Range("E38").Value = "  -2.02%  "

$ws.Range("D39").Value = "'3.41"
$ws.Range("E39").Value = "  -10.64%  "

$ws.Range("D40").Value = "'0.138"
$ws.Range("E40").Value = "  +0.99%  "

$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("D42").Value = "'5.79"
$ws.Range("E42").Value = "  -0.33%  "

$ws.Range("E43").Value = "  +0.37%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").Value = "'0.305"
$ws.Range("E45").Value = "  -1.65%  "

$ws.Range("D46").Value = "'1.95"
$ws.Range("E46").Value = "  +0.72%  "

$ws.Range("D47").Value = "'43.35"
$ws.Range("E47").Value = "  +11.45%  "

$ws.Range("D48").Value = "'8.59"
$ws.Range("E48").Value = "  -1.08%  "

$ws.Range("D49").Value = "'45.79"
$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("D50").Value = "'393.41"
$ws.Range("E50").Value = "  -0.68%  "

$ws.Range("D51").Value = "'145.14"
$ws.Range("E51").Value = "  +5.20%  "

# Clear the "number stored as text" quote-prefix marker the apostrophe entry added,
# restoring each forced cell to the workbook default style (matches source formatting).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
